# Update "想去人数" (want-to-go count) figures in the two sheets that carry
# the full per-event table: "展览" (sheet1) and "全部类型" (sheet4).
# Sheet "演出" and "本地生活" are unaffected by this refresh.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览" = @{
        "F2"  = 13863
        "F3"  = 92
        "F7"  = 1210
        "F8"  = 1035
        "F9"  = 13906
        "F10" = 14816
        "F14" = 177
        "F20" = 21
        "F22" = 19
        "F23" = 1149
        "F26" = 5757
        "F28" = 1061
        "F29" = 5422
        "F30" = 49
        "F31" = 52
        "F32" = 281
    }
    "全部类型" = @{
        "F2"  = 13863
        "F3"  = 92
        "F8"  = 1210
        "F9"  = 1035
        "F10" = 13906
        "F11" = 14816
        "F15" = 177
        "F21" = 21
        "F23" = 19
        "F24" = 1149
        "F27" = 5757
        "F29" = 1061
        "F30" = 5422
        "F31" = 49
        "F32" = 52
        "F33" = 281
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellMap = $updates[$sheetName]
    foreach ($addr in $cellMap.Keys) {
        $ws.Range($addr).Value = $cellMap[$addr]
    }
}
